$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note text ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$text = @"
Conversión del día 💰
✅ Dólar paralelo: 68

Binance
✅ 1000 Bs = 4.42 = 17426.3 pesos
✅ 17426.3 pesos = 4.39 = 932.36 Bs

Promedio competencia
✅ Tasa pesos: 20
✅ Tasa Bs: 20
✅ % Ganancia: 20%
"@
$ws1.Range("A1").Value = $text

# --- tasas: update the rate table values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 226.402
$ws2.Range("O10").Value = 3945.35
$ws2.Range("N12").Value = 3968
$ws2.Range("O12").Value = 212.3
